# "Schedule" sheet: rows 9 and 11 of column A held the shared formula
# (=<prev row>+7) that generates the weekly date sequence. The author
# overtyped those two cells with fixed literal dates, so they stop being
# formulas (and fall out of the shared-formula group) while every other
# row in the sequence keeps computing normally.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Row 9 (A9) used to compute 46053 (01/31/2026); it's now pinned to
# 02/01/2026 as a plain value.
$ws.Range("A9").Value = 46054

# Row 11 (A11) is pinned to its previously-computed date, 02/07/2026,
# as a plain value instead of a formula.
$ws.Range("A11").Value = 46060
